$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3486.375
$ws.Range("I28").Value = 843.7778
$ws.Range("K28").Value = 843.7778
$ws.Range("M28").Value = -358.7778
$ws.Range("H38").Value = 3571.889
$ws.Range("I38").Value = 1863.091
$ws.Range("J38").Value = 6257.143
$ws.Range("K38").Value = 5589.272999999999
$ws.Range("L38").Value = 18771.429
$ws.Range("M38").Value = -5217.272999999999
$ws.Range("N38").Value = -19515.429
$ws.Range("H69").Value = 6317.9
$ws.Range("I69").Value = 2671
$ws.Range("J69").Value = 6961.4707
$ws.Range("K69").Value = 8013
$ws.Range("L69").Value = 20884.4121
$ws.Range("M69").Value = -7139
$ws.Range("N69").Value = -22632.4121
$ws.Range("H72").Value = 6317.9
$ws.Range("I72").Value = 2671
$ws.Range("J72").Value = 6961.4707
$ws.Range("K72").Value = 24039
$ws.Range("L72").Value = 62653.2363
$ws.Range("M72").Value = -19671
$ws.Range("N72").Value = -71389.23629999999
$ws.Range("H86").Value = 5040.6
$ws.Range("I86").Value = 4050.75
$ws.Range("J86").Value = 9000
$ws.Range("K86").Value = 4050.75
$ws.Range("L86").Value = 9000
$ws.Range("M86").Value = -2927.75
$ws.Range("N86").Value = -11246
$ws.Range("H89").Value = 5040.6
$ws.Range("I89").Value = 4050.75
$ws.Range("J89").Value = 9000
$ws.Range("K89").Value = 20253.75
$ws.Range("L89").Value = 45000
$ws.Range("M89").Value = -14637.75
$ws.Range("N89").Value = -56232
$ws.Range("H113").Value = 3499.5
$ws.Range("I113").Value = 3499.5
$ws.Range("K113").Value = 3499.5
$ws.Range("M113").Value = -245.5
$ws.Range("H137").Value = 2632.3655
$ws.Range("I137").Value = 1664.2963
$ws.Range("J137").Value = 3677.88
$ws.Range("K137").Value = 4992.8889
$ws.Range("L137").Value = 11033.64
$ws.Range("M137").Value = -2442.8889
$ws.Range("N137").Value = -16133.64
$ws.Range("H140").Value = 50780
$ws.Range("J140").Value = 50780
$ws.Range("L140").Value = 50780
$ws.Range("N140").Value = -61140

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 204.64706
$ws.Range("J4").Value = 291.6
$ws.Range("L4").Value = 291.6
$ws.Range("N4").Value = -523.6
$ws.Range("H5").Value = 81.25
$ws.Range("I5").Value = 25.5
$ws.Range("K5").Value = 25.5
$ws.Range("M5").Value = 86.5
$ws.Range("H45").Value = 2597
$ws.Range("I45").Value = 2331.1667
$ws.Range("J45").Value = 2862.8333
$ws.Range("K45").Value = 2331.1667
$ws.Range("L45").Value = 2862.8333
$ws.Range("M45").Value = -1954.1667
$ws.Range("N45").Value = -3616.8333
$ws.Range("H74").Value = 4923
$ws.Range("I74").Value = 4923
$ws.Range("K74").Value = 4923
$ws.Range("M74").Value = -4049
$ws.Range("H77").Value = 4923
$ws.Range("I77").Value = 4923
$ws.Range("K77").Value = 24615
$ws.Range("M77").Value = -20247
$ws.Range("H114").Value = 7512000
$ws.Range("J114").Value = 7512000
$ws.Range("L114").Value = 7512000
$ws.Range("N114").Value = -7520678
$ws.Range("H122").Value = 4600
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4600
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 13800
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -18700

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 81.25
$ws.Range("I4").Value = 25.5
$ws.Range("K4").Value = 25.5
$ws.Range("M4").Value = 89.5
$ws.Range("H20").Value = 3007.5
$ws.Range("I20").Value = 3007.5
$ws.Range("K20").Value = 3007.5
$ws.Range("M20").Value = -2760.5
$ws.Range("H99").Value = 2151.8125
$ws.Range("I99").Value = 2172.2222
$ws.Range("J99").Value = 2125.5715
$ws.Range("K99").Value = 2172.2222
$ws.Range("L99").Value = 2125.5715
$ws.Range("M99").Value = -674.2222000000002
$ws.Range("N99").Value = -5121.5715

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H58").Value = 3017.125
$ws.Range("I58").Value = 2551
$ws.Range("J58").Value = 4042.6
$ws.Range("K58").Value = 2551
$ws.Range("L58").Value = 4042.6
$ws.Range("M58").Value = -2348
$ws.Range("N58").Value = -4448.6
$ws.Range("H136").Value = 3017.125
$ws.Range("I136").Value = 2551
$ws.Range("J136").Value = 4042.6
$ws.Range("K136").Value = 7653
$ws.Range("L136").Value = 12127.8
$ws.Range("M136").Value = -5103
$ws.Range("N136").Value = -17227.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 74.60869599999999
$ws.Range("I2").Value = 183.85715
$ws.Range("J2").Value = 26.8125
$ws.Range("K2").Value = 1103.1429
$ws.Range("L2").Value = 160.875
$ws.Range("M2").Value = -990.1428999999998
$ws.Range("N2").Value = -386.875
$ws.Range("H68").Value = 1885
$ws.Range("I68").Value = 1900
$ws.Range("J68").Value = 1881.25
$ws.Range("K68").Value = 5700
$ws.Range("L68").Value = 5643.75
$ws.Range("M68").Value = -4889
$ws.Range("N68").Value = -7265.75
$ws.Range("H71").Value = 1885
$ws.Range("I71").Value = 1900
$ws.Range("J71").Value = 1881.25
$ws.Range("K71").Value = 17100
$ws.Range("L71").Value = 16931.25
$ws.Range("M71").Value = -13044
$ws.Range("N71").Value = -25043.25
$ws.Range("H103").Value = 487.14285
$ws.Range("J103").Value = 783.3333
$ws.Range("L103").Value = 2349.9999
$ws.Range("N103").Value = -4107.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 204.57895
$ws.Range("I2").Value = 81.90000000000001
$ws.Range("J2").Value = 340.8889
$ws.Range("K2").Value = 81.90000000000001
$ws.Range("L2").Value = 340.8889
$ws.Range("M2").Value = 31.09999999999999
$ws.Range("N2").Value = -566.8888999999999
$ws.Range("H3").Value = 20777982
$ws.Range("I3").Value = 20864426
$ws.Range("J3").Value = 20000000
$ws.Range("K3").Value = 20864426
$ws.Range("L3").Value = 20000000
$ws.Range("M3").Value = -20864310
$ws.Range("N3").Value = -20000232
$ws.Range("H11").Value = 11977100
$ws.Range("I11").Value = 11346375
$ws.Range("K11").Value = 11346375
$ws.Range("M11").Value = -11346236
$ws.Range("H26").Value = 34166.5
$ws.Range("J26").Value = 34166.5
$ws.Range("L26").Value = 34166.5
$ws.Range("N26").Value = -34726.5
$ws.Range("H50").Value = 34166.5
$ws.Range("J50").Value = 34166.5
$ws.Range("L50").Value = 34166.5
$ws.Range("N50").Value = -35162.5
$ws.Range("H80").Value = 3155.4285
$ws.Range("I80").Value = 2940
$ws.Range("J80").Value = 3694
$ws.Range("K80").Value = 2940
$ws.Range("L80").Value = 3694
$ws.Range("M80").Value = -1942
$ws.Range("N80").Value = -5690
$ws.Range("H83").Value = 3155.4285
$ws.Range("I83").Value = 2940
$ws.Range("J83").Value = 3694
$ws.Range("K83").Value = 14700
$ws.Range("L83").Value = 18470
$ws.Range("M83").Value = -9708
$ws.Range("N83").Value = -28454
$ws.Range("H102").Value = 2174.923
$ws.Range("I102").Value = 1997.6364
$ws.Range("J102").Value = 3150
$ws.Range("K102").Value = 1997.6364
$ws.Range("L102").Value = 3150
$ws.Range("M102").Value = -375.6364000000001
$ws.Range("N102").Value = -6394
$ws.Range("H107").Value = 1087.7222
$ws.Range("J107").Value = 1218
$ws.Range("L107").Value = 1218
$ws.Range("N107").Value = -5058
$ws.Range("H122").Value = 2582.4
$ws.Range("I122").Value = 2228
$ws.Range("K122").Value = 6684
$ws.Range("M122").Value = -4234
$ws.Range("H126").Value = 3086
$ws.Range("I126").Value = 3086
$ws.Range("K126").Value = 9258
$ws.Range("M126").Value = -6788

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2690
$ws.Range("I40").Value = 2690
$ws.Range("K40").Value = 2690
$ws.Range("M40").Value = -2554
$ws.Range("H46").Value = 5828.8125
$ws.Range("J46").Value = 8793.799999999999
$ws.Range("L46").Value = 8793.799999999999
$ws.Range("N46").Value = -9169.799999999999
$ws.Range("H82").Value = 6742.857
$ws.Range("J82").Value = 6742.857
$ws.Range("L82").Value = 6742.857
$ws.Range("N82").Value = -7464.857
$ws.Range("H85").Value = 6742.857
$ws.Range("J85").Value = 6742.857
$ws.Range("L85").Value = 6742.857
$ws.Range("N85").Value = -9238.857
$ws.Range("H132").Value = 3825
$ws.Range("I132").Value = 3766.6667
$ws.Range("K132").Value = 11300.0001
$ws.Range("M132").Value = -8770.000100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 73332.664
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H126").Value = 6734.2
$ws.Range("I126").Value = 3780.6667
$ws.Range("K126").Value = 11342.0001
$ws.Range("M126").Value = -8872.000100000001
